$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "User" entity block (rows 22-26) and the trailing comment row (31)
# down by 5 rows, opening up space for the new "Phylum" entity block.
$ws.Rows("22:26").Insert()

# New field in the Species block: a ManyToOne "phylum" relation.
$ws.Range("A21").Value = "phylum"
$ws.Range("B21").Value = "relation"
$ws.Range("D21").Value = "Species"
$ws.Range("E21").Value = "ManyToOne"

# New "Phylum" entity block.
$ws.Range("A23").Value = "phylumNameWorms"
$ws.Range("B23").Value = "text"
$ws.Range("D23").Value = "Phylum"

$ws.Range("A24").Value = "species"
$ws.Range("B24").Value = "relation"
$ws.Range("D24").Value = "Phylum"
$ws.Range("E24").Value = "OneToMany"

$ws.Range("D24").Select()
